$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "August (through 08-30)"
$ws.Range("B9").Value = 31
$ws.Range("C9").Value = 78
$ws.Range("D9").Value = 86
$ws.Range("E9").Value = 65
$ws.Range("G9").Value = 161
$ws.Range("H9").Value = 152

$ws.Range("B10").Value = 193
$ws.Range("C10").Value = 380
$ws.Range("D10").Value = 551
$ws.Range("E10").Value = 490
$ws.Range("G10").Value = 782
$ws.Range("H10").Value = 1065
